$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 0.12315217332237
$ws.Range("J4").Value = 0.4834173371618073
$ws.Range("K4").Value = 0.4079614427643469
$ws.Range("L4").Value = 2.682563357569987
